$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value2 = '30.395.08'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value2 = '  -0.46%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value2 = '1.921.33'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value2 = '  +3.60%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value2 = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value2 = '  +0.04%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '240.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value2 = '  +2.81%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value2 = '  +0.03%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '0.4740'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value2 = '  -0.08%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('B8').Value2 = 'Cardano'
$ws.Range('C8').Value2 = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.2846'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value2 = '  +3.72%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('B9').Value2 = 'Dogecoin'
$ws.Range('C9').Value2 = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.06580'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value2 = '  +4.28%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('B10').Value2 = 'Solana'
$ws.Range('C10').Value2 = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '19.07'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value2 = '  +7.77%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('B11').Value2 = 'Litecoin'
$ws.Range('C11').Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '104.79'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value2 = '  +23.89%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('B12').Value2 = 'WrappedEther'
$ws.Range('C12').Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '1.910.56'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value2 = '  +3.19%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').Value2 = 'TRON'
$ws.Range('C13').Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '0.07584'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value2 = '  +1.96%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('B14').Value2 = 'Polkadot'
$ws.Range('C14').Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '5.117'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value2 = '  +2.79%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('B15').Value2 = 'Polygon'
$ws.Range('C15').Value2 = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '0.6530'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value2 = '  +4.26%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').Value2 = 'BitcoinCash'
$ws.Range('C16').Value2 = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '299.74'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value2 = '  +21.33%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value2 = 'WrappedBTC'
$ws.Range('C17').Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '30.418.15'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value2 = '  -0.22%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('B18').Value2 = 'Dai'
$ws.Range('C18').Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '1.001'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value2 = '  +0.03%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('B19').Value2 = 'Avalanche'
$ws.Range('C19').Value2 = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '12.90'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value2 = '  +1.79%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('B20').Value2 = 'ShibaInu'
$ws.Range('C20').Value2 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '0.000007509'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value2 = '  +2.53%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('B21').Value2 = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '2.162.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value2 = '  +2.85%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('B22').Value2 = 'BinanceUSD'
$ws.Range('C22').Value2 = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '1.002'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value2 = '  +0.13%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('B23').Value2 = 'Uniswap'
$ws.Range('C23').Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '5.223'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value2 = '  +5.81%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('B24').Value2 = 'Chainlink'
$ws.Range('C24').Value2 = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '6.274'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value2 = '  +6.13%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('B25').Value2 = 'Monero'
$ws.Range('C25').Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '166.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value2 = '  +2.50%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('B26').Value2 = 'Cosmos'
$ws.Range('C26').Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '9.171'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value2 = '  +0.56%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('B27').Value2 = 'EthereumClassic'
$ws.Range('C27').Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '19.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value2 = '  +9.35%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('B28').Value2 = 'LidoDAOToken'
$ws.Range('C28').Value2 = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '2.014'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value2 = '  +7.46%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('B29').Value2 = 'Stellar'
$ws.Range('C29').Value2 = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '0.1118'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value2 = '  +9.61%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('B30').Value2 = 'Toncoin'
$ws.Range('C30').Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '1.355'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value2 = '  +0.19%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Value2 = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '4.097'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value2 = '  +1.93%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value2 = 'Filecoin'
$ws.Range('C32').Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '3.914'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value2 = '  +2.19%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value2 = 'Hedera'
$ws.Range('C33').Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '0.05004'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value2 = '  +3.43%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('B34').Value2 = 'ImmutableX'
$ws.Range('C34').Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '0.7374'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value2 = '  +5.24%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('B35').Value2 = 'ARBITRUM'
$ws.Range('C35').Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '1.140'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value2 = '  +0.60%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('B36').Value2 = 'Frax'
$ws.Range('C36').Value2 = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '0.9997'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value2 = '  +0.03%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '2.716'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '0.01946'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value2 = '  +2.60%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '2.695'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value2 = '  +0.28%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '2.043'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value2 = '  +2.26%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '0.8724'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value2 = '  -0.24%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '107.07'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value2 = '  +0.55%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '5.778'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value2 = '  +4.33%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '0.9994'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value2 = '  -0.09%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '69.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value2 = '  +9.96%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '0.4112'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value2 = '  +1.61%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '7.190'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value2 = '  -0.09%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '9.189'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value2 = '  +7.38%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value2 = 'Elrond'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '34.58'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value2 = '  +3.08%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value2 = 'Algorand'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '0.1201'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value2 = '  -0.16%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '0.05620'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value2 = '  +1.65%  '
$ws.Range('E51').Style = 'Normal'
